# Apply the "Generate Report for Handback" edit described by the diff.
#
# Summary of changes:
#  - Status text everywhere it reads "Ready for handoff" becomes
#    "Handed back: in sync with en-US" (Overview!E2:F2, zh-cn!C2, de-de!C2).
#  - zh-cn / de-de sheets get their "Latest Target File" (I2) and
#    "Latest Handback File" (J2) populated, I2 becomes a hyperlink to the
#    source markdown file (same look as A2), and "Latest Handback DateTime"
#    (K2) gets a real timestamp instead of the epoch placeholder.
#  - Column widths for the columns whose text got longer are widened.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdFileName = "9a9450c5-c191-4572-8715-063a4d1820e2.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9861a96521ad56b3921a39fd5ff02e796453eeb7/e2e/9a9450c5-c191-4572-8715-063a4d1820e2.md"

# ---------------------------------------------------------------------
# 1. Status -> "Handed back: in sync with en-US" everywhere it shows up.
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws3.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn row 2: Latest Target File / Latest Handback File / Latest
#    Handback DateTime.
# ---------------------------------------------------------------------
$ws2.Range("I2").Value = $mdFileName
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null
$ws2.Range("I2").Font.Underline = 2
$ws2.Range("I2").Font.Color = 15570276

$ws2.Range("J2").Value = "9a9450c5-c191-4572-8715-063a4d1820e2.8cf3e283a5c6ba46ce8bfeea422d9646d4919405.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-02 05:03:34"

# ---------------------------------------------------------------------
# 3. de-de row 2: Latest Target File / Latest Handback File / Latest
#    Handback DateTime.
# ---------------------------------------------------------------------
$ws3.Range("I2").Value = $mdFileName
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null
$ws3.Range("I2").Font.Underline = 2
$ws3.Range("I2").Font.Color = 15570276

$ws3.Range("J2").Value = "9a9450c5-c191-4572-8715-063a4d1820e2.8cf3e283a5c6ba46ce8bfeea422d9646d4919405.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-02 05:03:42"

# ---------------------------------------------------------------------
# 4. Column widths: widen columns whose cell text got longer.
# ---------------------------------------------------------------------
$ws1.Columns("E").ColumnWidth = 29.14
$ws1.Columns("F").ColumnWidth = 29.14

$ws2.Columns("C").ColumnWidth = 29.14
$ws2.Columns("I").ColumnWidth = 39.17
$ws2.Columns("J").ColumnWidth = 39.17

$ws3.Columns("C").ColumnWidth = 29.14
$ws3.Columns("I").ColumnWidth = 39.17
$ws3.Columns("J").ColumnWidth = 39.17
